$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7595
$ws.Range("J3").Value = 7989
$ws.Range("J4").Value = 1736
$ws.Range("J5").Value = 622
$ws.Range("J6").Value = 10918
$ws.Range("J7").Value = 28860

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 282
$ws.Range("J7").Value = 434

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 482
$ws.Range("J3").Value = 522
$ws.Range("J4").Value = 94
$ws.Range("J6").Value = 683
$ws.Range("J7").Value = 1828

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 213
$ws.Range("J6").Value = 159
$ws.Range("J7").Value = 586

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J6").Value = 466
$ws.Range("J7").Value = 1302

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 270
$ws.Range("J6").Value = 257
$ws.Range("J7").Value = 888

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 210
$ws.Range("J5").Value = 20
$ws.Range("J7").Value = 724

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 173
$ws.Range("J7").Value = 439

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 229
$ws.Range("J8").Value = 1828
$ws.Range("J11").Value = 526
$ws.Range("J15").Value = 354
$ws.Range("J17").Value = 38
$ws.Range("J18").Value = 232
$ws.Range("J19").Value = 841
$ws.Range("J23").Value = 265
$ws.Range("J27").Value = 177
$ws.Range("J29").Value = 1539
$ws.Range("J31").Value = 305
$ws.Range("J33").Value = 1302
$ws.Range("J36").Value = 393
$ws.Range("J37").Value = 888
$ws.Range("J42").Value = 1224
$ws.Range("J44").Value = 227
$ws.Range("J47").Value = 208
$ws.Range("J50").Value = 180
$ws.Range("J51").Value = 363
$ws.Range("J52").Value = 734
$ws.Range("J53").Value = 434
$ws.Range("J55").Value = 453
$ws.Range("J58").Value = 17
$ws.Range("J59").Value = 35
$ws.Range("J60").Value = 172
$ws.Range("J63").Value = 88
$ws.Range("J65").Value = 724
$ws.Range("J67").Value = 1048
$ws.Range("J78").Value = 338
$ws.Range("J79").Value = 792
$ws.Range("J83").Value = 586
$ws.Range("J85").Value = 1181
$ws.Range("J86").Value = 174
$ws.Range("J89").Value = 360
$ws.Range("J91").Value = 332
$ws.Range("J93").Value = 121
$ws.Range("J94").Value = 319
$ws.Range("J99").Value = 439
$ws.Range("J101").Value = 28860

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 99
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 267
$ws.Range("J6").Value = 291
$ws.Range("J7").Value = 1048

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 467
$ws.Range("J6").Value = 392
$ws.Range("J7").Value = 1539

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 209
$ws.Range("J6").Value = 328
$ws.Range("J7").Value = 841

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 68
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 227

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 247
$ws.Range("J6").Value = 648
$ws.Range("J7").Value = 1224

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 256
$ws.Range("J7").Value = 453

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 89
$ws.Range("J7").Value = 265

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 90
$ws.Range("J3").Value = 134
$ws.Range("J7").Value = 332

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 266
$ws.Range("J7").Value = 792

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 128
$ws.Range("J7").Value = 393

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 171
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 73
$ws.Range("J7").Value = 354

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 143
$ws.Range("J7").Value = 526

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 108
$ws.Range("J7").Value = 360

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 47
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 177

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 174

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 77
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 363

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 59
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 429
$ws.Range("J7").Value = 1181

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 206
$ws.Range("J7").Value = 734

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 17
